$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple +1 increments on column B (Taxonsorteringsordning) ---
$ws.Range("B2").Value = 79864
$ws.Range("B3").Value = 79245
$ws.Range("B4").Value = 79864
$ws.Range("B5").Value = 79002
$ws.Range("B9").Value = 79864
$ws.Range("B10").Value = 79245
$ws.Range("B11").Value = 79245
$ws.Range("B12").Value = 79245
$ws.Range("B13").Value = 79245
$ws.Range("B15").Value = 79835
$ws.Range("B16").Value = 79245
$ws.Range("B17").Value = 79245
$ws.Range("B21").Value = 79835
$ws.Range("B24").Value = 79835
$ws.Range("B25").Value = 79003
$ws.Range("B30").Value = 79245
$ws.Range("B31").Value = 79245

# --- Row 6/7/8: content rotates (row6->row7, row7->row8, row8->row6), ---
# --- with the Taxonsorteringsordning (col B) value also getting the   ---
# --- same +1 remap applied to the moved-in record.                    ---
$ws.Range("A6").Value = 131198844
$ws.Range("B6").Value = 79245
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("M6").Value = ""
$ws.Range("Q6").Value = 466309
$ws.Range("R6").Value = 6789077

$ws.Range("A7").Value = 131197802
$ws.Range("B7").Value = 57881
$ws.Range("E7").Value = 100049
$ws.Range("F7").Value = "Spillkråka"
$ws.Range("G7").Value = "Dryocopus martius"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("M7").Value = "äldre spår"
$ws.Range("Q7").Value = 465938
$ws.Range("R7").Value = 6789021

$ws.Range("A8").Value = 131198466
$ws.Range("B8").Value = 79835
$ws.Range("E8").Value = 229821
$ws.Range("F8").Value = "Vedflamlav"
$ws.Range("G8").Value = "Ramboldia elabens"
$ws.Range("H8").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q8").Value = 466092
$ws.Range("R8").Value = 6789074

# --- Row 28/29: whole rows swap content ---
$ws.Range("A28").Value = 131198231
$ws.Range("B28").Value = 79245
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("M28").Value = ""
$ws.Range("Q28").Value = 466052
$ws.Range("R28").Value = 6789006
$ws.Range("S28").Value = 50
$ws.Range("AC28").Value = "Rikligt till måttligt i en radie av ca 50 meter  1 bild gren i förgrund"

$ws.Range("A29").Value = 131199091
$ws.Range("B29").Value = 57884
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = "Tretåig hackspett"
$ws.Range("G29").Value = "Picoides tridactylus"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("M29").Value = "färska spår"
$ws.Range("Q29").Value = 466114
$ws.Range("R29").Value = 6788962
$ws.Range("S29").Value = 10
$ws.Range("AC29").Value = "2 bild, tall med gran till vänster."
